# Add cantrals by cantons
# Restructure the SH 2015 sheet: collapse the old two-row header (units
# split across two rows, with a stray partial data row in between) into a
# single header row with new column titles (idx, idx2, Name, Date Start,
# Date End, (m3/s), (MW1), (MW2), (GWh) Winter, (GWh) Summer, (GWh) Year),
# and drop the orphan partial row that had no A:E values.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- 1. Remove the stray "continuation" row (old row 5: F=3, G=5.3, H=5.3,
#        no A:E values) and the old units row (old row 2) -----------------
# Delete the orphan row first (row 5) so row numbers above it are untouched
# while we still need them, then delete the old units row (row 2).
$ws.Rows.Item(5).Delete() | Out-Null
$ws.Rows.Item(2).Delete() | Out-Null

# After these two deletions the sheet now looks like:
#   Row 1: old first header row (French labels) - about to be overwritten
#   Row 2: Wunderklingen data (was row 3)
#   Row 3: Engeweiher data   (was row 4)
#   Row 4: Eglisau data      (was row 6)
#   Row 5: Neuhausen data    (was row 7)
#   Row 6: Rheinau data      (was row 8)
#   Row 7: Schaffhausen data (was row 9)

# --- 2. Rewrite the header row -------------------------------------------
$ws.Range("A1").Value = "idx"
$ws.Range("B1").Value = "idx2"
$ws.Range("C1").Value = "Name"
$ws.Range("D1").Value = "Date Start"
$ws.Range("E1").Value = "Date End"
$ws.Range("F1").Value = "(m3/s)"
$ws.Range("G1").Value = "(MW1)"
$ws.Range("H1").Value = "(MW2)"
$ws.Range("I1").Value = "(GWh) Winter"
$ws.Range("J1").Value = "(GWh) Summer"
$ws.Range("K1").Value = "(GWh) Year"

# A1:E1 keep the plain default style (no explicit style). F1:K1 get a new
# "font-only" header style (Arial 9, no explicit number format) - build it
# as a transient named style so the resulting xf has no applyNumberFormat
# flag, then drop the named style again so only the xf record remains.
$ws.Range("A1:E1").Style = "Normal"

$headerStyle = $wb.Styles.Add("SH2015Header")
$headerStyle.Font.Bold = $true
$headerStyle.Font.Size = 9
$headerStyle.Font.Bold = $false
$ws.Range("F1:K1").Style = "SH2015Header"
$wb.Styles.Item("SH2015Header").Delete()

# --- 3. Fix up the sheet view / selection / used range --------------------
$ws.Range("A4:K4").Select()

$wb.Save()
